$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H133").Value = 93954.42999999999
$ws.Range("J133").Value = 93954.42999999999
$ws.Range("L133").Value = 93954.42999999999
$ws.Range("N133").Value = -104074.43
$ws.Range("H137").Value = 8775915
$ws.Range("I137").Value = 752.61536
$ws.Range("J137").Value = 13338999
$ws.Range("K137").Value = 2257.84608
$ws.Range("L137").Value = 40016997
$ws.Range("M137").Value = 292.1539199999997
$ws.Range("N137").Value = -40022097
$ws.Range("H140").Value = 59583.8
$ws.Range("J140").Value = 58347.668
$ws.Range("L140").Value = 58347.668
$ws.Range("N140").Value = -68707.66800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21500.75
$ws.Range("I61").Value = 39627.75
$ws.Range("J61").Value = 3373.75
$ws.Range("K61").Value = 39627.75
$ws.Range("L61").Value = 3373.75
$ws.Range("M61").Value = -39415.75
$ws.Range("N61").Value = -3797.75
$ws.Range("H74").Value = 8930156
$ws.Range("I74").Value = 12501125
$ws.Range("K74").Value = 12501125
$ws.Range("M74").Value = -12500251
$ws.Range("H77").Value = 8930156
$ws.Range("I77").Value = 12501125
$ws.Range("K77").Value = 62505625
$ws.Range("M77").Value = -62501257
$ws.Range("H136").Value = 21500.75
$ws.Range("I136").Value = 39627.75
$ws.Range("J136").Value = 3373.75
$ws.Range("K136").Value = 118883.25
$ws.Range("L136").Value = 10121.25
$ws.Range("M136").Value = -116333.25
$ws.Range("N136").Value = -15221.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2803.276
$ws.Range("I20").Value = 2359.5454
$ws.Range("J20").Value = 4197.857
$ws.Range("K20").Value = 2359.5454
$ws.Range("L20").Value = 4197.857
$ws.Range("M20").Value = -2112.5454
$ws.Range("N20").Value = -4691.857
$ws.Range("H22").Value = 527
$ws.Range("I22").Value = 489.7
$ws.Range("K22").Value = 489.7
$ws.Range("M22").Value = -316.7
$ws.Range("H94").Value = 508147.84
$ws.Range("I94").Value = 623257.9399999999
$ws.Range("J94").Value = 1663.6
$ws.Range("K94").Value = 623257.9399999999
$ws.Range("L94").Value = 1663.6
$ws.Range("M94").Value = -622806.9399999999
$ws.Range("N94").Value = -2565.6
$ws.Range("H105").Value = 2826.9412
$ws.Range("I105").Value = 2628.9167
$ws.Range("J105").Value = 3302.2
$ws.Range("K105").Value = 2628.9167
$ws.Range("L105").Value = 3302.2
$ws.Range("M105").Value = -881.9167000000002
$ws.Range("N105").Value = -6796.2
$ws.Range("H107").Value = 1465
$ws.Range("I107").Value = 1439.0526
$ws.Range("K107").Value = 1439.0526
$ws.Range("M107").Value = 480.9474
$ws.Range("H134").Value = 4808.143
$ws.Range("I134").Value = 1925
$ws.Range("K134").Value = 5775
$ws.Range("M134").Value = -3240

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 888.8889
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1183.3334
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1183.3334
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -1883.3334
$ws.Range("H31").Value = 4891.4883
$ws.Range("I31").Value = 1648.4584
$ws.Range("K31").Value = 1648.4584
$ws.Range("M31").Value = -1353.4584
$ws.Range("H34").Value = 4891.4883
$ws.Range("I34").Value = 1648.4584
$ws.Range("K34").Value = 1648.4584
$ws.Range("M34").Value = -1446.4584
$ws.Range("H58").Value = 504129.8
$ws.Range("J58").Value = 4622.25
$ws.Range("L58").Value = 4622.25
$ws.Range("N58").Value = -5028.25
$ws.Range("H86").Value = 7076.8
$ws.Range("I86").Value = 6128
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 6128
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -5005
$ws.Range("N86").Value = -10746
$ws.Range("H89").Value = 7076.8
$ws.Range("I89").Value = 6128
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 30640
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -25024
$ws.Range("N89").Value = -53732
$ws.Range("H132").Value = 11918481
$ws.Range("I132").Value = 13344099
$ws.Range("J132").Value = 38333
$ws.Range("K132").Value = 40032297
$ws.Range("L132").Value = 114999
$ws.Range("M132").Value = -40029767
$ws.Range("N132").Value = -120059
$ws.Range("H134").Value = 3922.0908
$ws.Range("J134").Value = 3999.3333
$ws.Range("L134").Value = 11997.9999
$ws.Range("N134").Value = -17067.9999
$ws.Range("H136").Value = 504129.8
$ws.Range("J136").Value = 4622.25
$ws.Range("L136").Value = 13866.75
$ws.Range("N136").Value = -18966.75
$ws.Range("H141").Value = 92193.53999999999
$ws.Range("I141").Value = 39899
$ws.Range("J141").Value = 96551.414
$ws.Range("K141").Value = 39899
$ws.Range("L141").Value = 96551.414
$ws.Range("M141").Value = -34719
$ws.Range("N141").Value = -106911.414

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1598.5385
$ws.Range("I22").Value = 2033.75
$ws.Range("J22").Value = 1405.1111
$ws.Range("K22").Value = 6101.25
$ws.Range("L22").Value = 4215.3333
$ws.Range("M22").Value = -5932.25
$ws.Range("N22").Value = -4553.3333
$ws.Range("H27").Value = 1598.5385
$ws.Range("I27").Value = 2033.75
$ws.Range("J27").Value = 1405.1111
$ws.Range("K27").Value = 6101.25
$ws.Range("L27").Value = 4215.3333
$ws.Range("M27").Value = -5999.25
$ws.Range("N27").Value = -4419.3333
$ws.Range("H33").Value = 424.75
$ws.Range("I33").Value = 383.16666
$ws.Range("J33").Value = 549.5
$ws.Range("K33").Value = 2298.99996
$ws.Range("L33").Value = 3297
$ws.Range("M33").Value = -2015.99996
$ws.Range("N33").Value = -3863
$ws.Range("H129").Value = 1959.5
$ws.Range("J129").Value = 1689.25
$ws.Range("L129").Value = 5067.75
$ws.Range("N129").Value = -15067.75
$ws.Range("H131").Value = 8336002.5
$ws.Range("J131").Value = 6538657.5
$ws.Range("L131").Value = 19615972.5
$ws.Range("N131").Value = -19626052.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 377.83334
$ws.Range("J97").Value = 384
$ws.Range("L97").Value = 384
$ws.Range("N97").Value = -1376
$ws.Range("H132").Value = 5527.524
$ws.Range("I132").Value = 5642.533
$ws.Range("J132").Value = 5240
$ws.Range("K132").Value = 16927.599
$ws.Range("L132").Value = 15720
$ws.Range("M132").Value = -14397.599
$ws.Range("N132").Value = -20780

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1210.1428
$ws.Range("I22").Value = 993.25
$ws.Range("K22").Value = 993.25
$ws.Range("M22").Value = -698.25
$ws.Range("H27").Value = 1210.1428
$ws.Range("I27").Value = 993.25
$ws.Range("K27").Value = 993.25
$ws.Range("M27").Value = -886.25
$ws.Range("H40").Value = 33334600
$ws.Range("I40").Value = 1406.4445
$ws.Range("K40").Value = 1406.4445
$ws.Range("M40").Value = -1270.4445
$ws.Range("H61").Value = 4496
$ws.Range("I61").Value = 4395.8
$ws.Range("K61").Value = 4395.8
$ws.Range("M61").Value = -4193.8
$ws.Range("H68").Value = 3793760.2
$ws.Range("I68").Value = 7578857.5
$ws.Range("J68").Value = 8663.333000000001
$ws.Range("K68").Value = 7578857.5
$ws.Range("L68").Value = 8663.333000000001
$ws.Range("M68").Value = -7578108.5
$ws.Range("N68").Value = -10161.333
$ws.Range("H71").Value = 3793760.2
$ws.Range("I71").Value = 7578857.5
$ws.Range("J71").Value = 8663.333000000001
$ws.Range("K71").Value = 37894287.5
$ws.Range("L71").Value = 43316.665
$ws.Range("M71").Value = -37890543.5
$ws.Range("N71").Value = -50804.665
$ws.Range("H82").Value = 15625750
$ws.Range("I82").Value = 31250000
$ws.Range("K82").Value = 31250000
$ws.Range("M82").Value = -31249639
$ws.Range("H85").Value = 15625750
$ws.Range("I85").Value = 31250000
$ws.Range("K85").Value = 31250000
$ws.Range("M85").Value = -31248752
$ws.Range("H93").Value = 3220.111
$ws.Range("I93").Value = 603.4
$ws.Range("J93").Value = 6491
$ws.Range("K93").Value = 603.4
$ws.Range("L93").Value = 6491
$ws.Range("M93").Value = 644.6
$ws.Range("N93").Value = -8987
$ws.Range("H100").Value = 3277.2
$ws.Range("I100").Value = 3277.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3277.2
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2736.2
$ws.Range("N100").Value = ""
$ws.Range("H113").Value = 4496
$ws.Range("I113").Value = 4395.8
$ws.Range("K113").Value = 4395.8
$ws.Range("M113").Value = -2225.8
$ws.Range("H132").Value = 6848.44
$ws.Range("I132").Value = 4808.3335
$ws.Range("K132").Value = 14425.0005
$ws.Range("M132").Value = -11895.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 300000
$ws.Range("J15").Value = 300000
$ws.Range("L15").Value = 300000
$ws.Range("N15").Value = -300576
$ws.Range("H96").Value = 9403.5
$ws.Range("J96").Value = 9403.5
$ws.Range("L96").Value = 9403.5
$ws.Range("N96").Value = -12149.5
$ws.Range("H107").Value = 2176.6296
$ws.Range("I107").Value = 2518.8572
$ws.Range("J107").Value = 978.8333
$ws.Range("K107").Value = 7556.571599999999
$ws.Range("L107").Value = 2936.4999
$ws.Range("M107").Value = -5636.571599999999
$ws.Range("N107").Value = -6776.4999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H132").Value = 26457062
$ws.Range("I132").Value = 3473658.5
$ws.Range("K132").Value = 10420975.5
$ws.Range("M132").Value = -10420975.5
